$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Baseline "plain" styling (style used by D5:J5 -- fillId=2, no border)
#    across the new block D18:J24, row by row. Specific cells below are then
#    overwritten with their own value + style.
# ---------------------------------------------------------------------------
$ws.Range("D5:J5").Copy()
$ws.Range("D18:J18").PasteSpecial(-4122)

$ws.Range("D5:J5").Copy()
$ws.Range("D19:J19").PasteSpecial(-4122)

$ws.Range("D5:J5").Copy()
$ws.Range("D20:J20").PasteSpecial(-4122)

$ws.Range("D5:J5").Copy()
$ws.Range("D21:J21").PasteSpecial(-4122)

$ws.Range("D5:J5").Copy()
$ws.Range("D22:J22").PasteSpecial(-4122)

$ws.Range("D5:J5").Copy()
$ws.Range("D23:J23").PasteSpecial(-4122)

$ws.Range("D5:J5").Copy()
$ws.Range("D24:J24").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Row 19 -- headers: "ε" | "t [s]" | "Períodos"
# ---------------------------------------------------------------------------
$ws.Range("E7").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$ws.Range("F19").Value = "t [s]"

$ws.Range("E7").Copy()
$ws.Range("G19").PasteSpecial(-4122)
$ws.Range("G19").Value = "Períodos"

$ws.Range("E7").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").Value = "ε"
$ws.Range("E19").Font.Name = "Calibri"
$ws.Range("E19").Font.Size = 11

# ---------------------------------------------------------------------------
# 3) Rows 20-23 -- data: epsilon (scientific) | t [s] | Períodos
# ---------------------------------------------------------------------------
# Row 20
$ws.Range("E7").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").NumberFormat = "0.00E+00"
$ws.Range("E20").Value = 0.01

$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "0,914"
$ws.Range("E7").Copy()
$ws.Range("F20").PasteSpecial(-4122)

$ws.Range("E7").Copy()
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("G20").Value = 22

# Row 21
$ws.Range("E7").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E21").NumberFormat = "0.00E+00"
$ws.Range("E21").Value = 0.001

$ws.Range("G8").Copy()
$ws.Range("F21").PasteSpecial(-4122)
$ws.Range("F21").Value = 1182

$ws.Range("E7").Copy()
$ws.Range("G21").PasteSpecial(-4122)
$ws.Range("G21").Value = 29

# Row 22
$ws.Range("E7").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "0.00E+00"
$ws.Range("E22").Value = 0.00001

$ws.Range("G8").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1718

$ws.Range("E7").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G22").Value = 42

# Row 23
$ws.Range("E7").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").NumberFormat = "0.00E+00"
$ws.Range("E23").Value = 0.0000000001

$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "3,06"
$ws.Range("E7").Copy()
$ws.Range("F23").PasteSpecial(-4122)

$ws.Range("E7").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("G23").Value = 76

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Column I width: 12.140625 custom -> 11.5703125 bestFit
# ---------------------------------------------------------------------------
$ws.Columns("I").ColumnWidth = 11.5703125

# ---------------------------------------------------------------------------
# 5) View state: scroll + selection
# ---------------------------------------------------------------------------
$ws.Range("G29").Select()
$excel.ActiveWindow.ScrollRow = 9

# ---------------------------------------------------------------------------
# 6) Leftover empty textbox (as present in the authored workbook)
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.AddTextbox(1, 468.0, 177.37496062992125, 0.005118110236220472, 13.561181102362205)
$shp.Name = "TextBox 1"
